# "Implemented Session Eviction API automation" / "Modified Authoring test cases"
#
#  - Existing rows 2-4 (S1_TC_T1..T3) get a new STATUS (col L) value of "PASS"
#  - Three new rows are appended for the Session-Eviction API tests:
#      row 5 = S1_TC_T4 "Evict user id"
#      row 6 = S1_TC_T5 "Activate evicted user"
#      row 7 = S1_TC_T6 "Get evicted user info or bucket information"
#  - Row 10 is touched only to extend the used range with a custom row height
#    (no cell data)
#  - Columns D (APIPATH) and H (BODY) are widened to fit the new, longer content
#  - The view's frozen "topLeftCell" scroll is dropped and the selection moves
#    to A8:XFD12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 2-4: add STATUS = PASS -----------------------------------
$ws.Cells.Item(2, 12).Value = "PASS"
$ws.Cells.Item(3, 12).Value = "PASS"
$ws.Cells.Item(4, 12).Value = "PASS"

# --- Row 5: S1_TC_T4 - "Evict user id" --------------------------------------
$ws.Cells.Item(5, 1).Value = "S1_TC_T4"
$ws.Cells.Item(5, 2).Value = "Evict user id"
$ws.Cells.Item(5, 3).Value = "1PAUTH"
$ws.Cells.Item(5, 4).Value = "/users/access"
$ws.Cells.Item(5, 4).Style = "Hyperlink"
$ws.Cells.Item(5, 5).Value = "POST"
$ws.Cells.Item(5, 6).Value = "Content-Type=application/json"
$ws.Cells.Item(5, 8).Value = '{"truid":"(SYS_USER1)" ,"userStatus":"Deactivate"}'
$ws.Cells.Item(5, 10).Value = "status=200||Success=OK"
$ws.Cells.Item(5, 12).Value = "PASS"

# --- Row 6: S1_TC_T5 - "Activate evicted user" ------------------------------
$ws.Cells.Item(6, 1).Value = "S1_TC_T5"
$ws.Cells.Item(6, 2).Value = "Activate evicted user"
$ws.Cells.Item(6, 3).Value = "1PAUTH"
$ws.Cells.Item(6, 4).Value = "/users/access"
$ws.Cells.Item(6, 4).Style = "Hyperlink"
$ws.Cells.Item(6, 5).Value = "POST"
$ws.Cells.Item(6, 6).Value = "Content-Type=application/json"
$ws.Cells.Item(6, 8).Value = '{"truid":"(SYS_USER1)" ,"userStatus":"Activate"}'
$ws.Cells.Item(6, 10).Value = "status=200||Success=OK"
$ws.Cells.Item(6, 12).Value = "PASS"

# --- Row 7: S1_TC_T6 - "Get evicted user info or bucket information" -------
$ws.Cells.Item(7, 1).Value = "S1_TC_T6"
$ws.Cells.Item(7, 2).Value = "Get evicted user info or bucket information"
$ws.Cells.Item(7, 3).Value = "1PAUTH"
$ws.Cells.Item(7, 4).Value = "/auth/bucket/f"
$ws.Cells.Item(7, 4).Style = "Hyperlink"
$ws.Cells.Item(7, 5).Value = "GET"
$ws.Cells.Item(7, 10).Value = "status=200||truid=(SYS_USER1)"
$ws.Cells.Item(7, 12).Value = "FAIL"
$ws.Rows.Item(7).RowHeight = 30

# --- Row 10: used range extends with a custom height, no data --------------
$ws.Rows.Item(10).RowHeight = 21

# --- Column widths: widen APIPATH (D) and BODY (H) for the new content -----
$ws.Columns.Item(4).ColumnWidth = 13.75
$ws.Columns.Item(8).ColumnWidth = 69.42

# --- View: drop the frozen scroll position, select A8:XFD12 ----------------
$ws.Range("A8:XFD12").Select()
